$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the dataset. It belongs at the
# top of the date-ordered block (row 14), so insert a new row there and push
# the existing rows 14-33 down to 15-34.
$ws.Rows.Item(14).Insert()

# Match the date-column style (s="2") used throughout the table by copying
# it down from the row above into the newly inserted row's date cell.
$ws.Cells.Item(13, 4).Copy()
$ws.Cells.Item(14, 4).PasteSpecial(-4122)

# Populate the new row 14 with the new observation's data.
$ws.Cells.Item(14, 1).Value2 = 8
$ws.Cells.Item(14, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(14, 3).Value2 = "Coquimbo"
$ws.Cells.Item(14, 4).Value2 = 44810
$ws.Cells.Item(14, 5).Value2 = 4
$ws.Cells.Item(14, 6).Value2 = "Fruta"
$ws.Cells.Item(14, 7).Value2 = 100101
$ws.Cells.Item(14, 8).Value2 = "Berries"
$ws.Cells.Item(14, 9).Value2 = 100101001
$ws.Cells.Item(14, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(14, 11).Value2 = "Sin especificar"
$ws.Cells.Item(14, 12).Value2 = "Primera"
$ws.Cells.Item(14, 13).Value2 = 200
$ws.Cells.Item(14, 14).Value2 = 15000
$ws.Cells.Item(14, 15).Value2 = 16000
$ws.Cells.Item(14, 16).Value2 = 15500
$ws.Cells.Item(14, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(14, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(14, 19).Value2 = 7750
$ws.Cells.Item(14, 20).Value2 = 2
